$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the scraped price/volume refresh from the crypto feed.
# Price cells are stored as literal text in the sheet (e.g. "7.230",
# "20.00"), so NumberFormat is forced to Text ("@") immediately before
# assigning any value that Excel would otherwise auto-convert/normalize
# into a number (dropping meaningful trailing zeros, multi-dot strings).

# Row 2
$ws.Range("D2").Value = "29.383.71"
$ws.Range("E2").Value = "  -0.17%  "

# Row 3
$ws.Range("D3").Value = "1.916.57"
$ws.Range("E3").Value = "  +0.85%  "

# Row 4
$ws.Range("E4").Value = "  +0.57%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.21"
$ws.Range("E5").Value = "  -0.49%  "

# Row 6
$ws.Range("E6").Value = "  +0.35%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4813"
$ws.Range("E7").Value = "  +0.22%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4061"
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08212"
$ws.Range("E9").Value = "  +1.81%  "

# Row 10
$ws.Range("E10").Value = "  +0.58%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.33"
$ws.Range("E11").Value = "  -0.47%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.076"
$ws.Range("E12").Value = "  +1.95%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.896.05"
$ws.Range("E13").Value = "  -0.66%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.230"
$ws.Range("E14").Value = "  +2.20%  "

# Row 15
$ws.Range("E15").Value = "  +1.70%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06856"
$ws.Range("E16").Value = "  +2.10%  "

# Row 17
$ws.Range("E17").Value = "  +0.43%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001038"
$ws.Range("E18").Value = "  +0.55%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.60"
$ws.Range("E19").Value = "  -0.15%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.007"

# Row 21
$ws.Range("D21").Value = "29.410.32"
$ws.Range("E21").Value = "  -0.13%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.655"
$ws.Range("E22").Value = "  +2.03%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.80"
$ws.Range("E23").Value = "  +0.32%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.183"
$ws.Range("E24").Value = "  +1.19%  "

# Row 25
$ws.Range("D25").Value = "2.139.04"
$ws.Range("E25").Value = "  -0.19%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.588"
$ws.Range("E26").Value = "  +8.02%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.74"
$ws.Range("E27").Value = "  +1.28%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.00"
$ws.Range("E28").Value = "  +0.79%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.111"
$ws.Range("E29").Value = "  +0.86%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.29"
$ws.Range("E30").Value = "  +1.50%  "

# Row 31
$ws.Range("E31").Value = "  -1.71%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09647"
$ws.Range("E32").Value = "  +1.60%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.624"
$ws.Range("E33").Value = "  +2.06%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.548"
$ws.Range("E34").Value = "  +0.04%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.374"
$ws.Range("E35").Value = "  -1.09%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02282"
$ws.Range("E36").Value = "  +1.40%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06096"
$ws.Range("E37").Value = "  +0.29%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.182"
$ws.Range("E38").Value = "  +0.72%  "

# Row 39
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.067"
$ws.Range("E39").Value = "  +2.01%  "

# Row 40
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.85"
$ws.Range("E40").Value = "  +6.25%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5948"
$ws.Range("E41").Value = "  +1.12%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1846"
$ws.Range("E42").Value = "  +0.08%  "

# Row 43
$ws.Range("E43").Value = "  -0.94%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.384"
$ws.Range("E44").Value = "  -0.22%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07592"
$ws.Range("E45").Value = "  -2.82%  "

# Row 46
$ws.Range("E46").Value = "  +1.12%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5579"
$ws.Range("E47").Value = "  +0.86%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.951"
$ws.Range("E48").Value = "  +1.49%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.70"
$ws.Range("E49").Value = "  +3.91%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.421"
$ws.Range("E50").Value = "  +3.28%  "

# Row 51
$ws.Range("E51").Value = "  -0.32%  "
